$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 333588.16
$ws.Range("I12").Value = 257.25
$ws.Range("K12").Value = 257.25
$ws.Range("M12").Value = -87.25
$ws.Range("H21").Value = 50019
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents() | Out-Null
$ws.Range("H23").Value = 50019
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents() | Out-Null
$ws.Range("H38").Value = 367.5
$ws.Range("I38").Value = 271.52942
$ws.Range("J38").Value = 1999
$ws.Range("K38").Value = 814.58826
$ws.Range("L38").Value = 5997
$ws.Range("M38").Value = -442.58826
$ws.Range("N38").Value = -6741
$ws.Range("H112").Value = 1401.303
$ws.Range("I112").Value = 718.4286
$ws.Range("J112").Value = 1453.2609
$ws.Range("K112").Value = 2155.2858
$ws.Range("L112").Value = 4359.7827
$ws.Range("M112").Value = -1047.2858
$ws.Range("N112").Value = -6575.7827
$ws.Range("H125").Value = 1223.2727
$ws.Range("I125").Value = 1172.5
$ws.Range("J125").Value = 1252.2858
$ws.Range("K125").Value = 10552.5
$ws.Range("L125").Value = 11270.5722
$ws.Range("M125").Value = -8092.5
$ws.Range("N125").Value = -16190.5722
$ws.Range("H129").Value = 16668029
$ws.Range("I129").Value = 83335450
$ws.Range("J129").Value = 1173.8334
$ws.Range("K129").Value = 250006350
$ws.Range("L129").Value = 3521.5002
$ws.Range("M129").Value = -250001350
$ws.Range("N129").Value = -13521.5002
$ws.Range("H137").Value = 2707136.2
$ws.Range("I137").Value = 4004900
$ws.Range("J137").Value = 3461.5833
$ws.Range("K137").Value = 12014700
$ws.Range("L137").Value = 10384.7499
$ws.Range("M137").Value = -12012150
$ws.Range("N137").Value = -15484.7499
$ws.Range("H138").Value = 6594.363
$ws.Range("I138").Value = 3033.8333
$ws.Range("J138").Value = 7869.776
$ws.Range("K138").Value = 9101.499899999999
$ws.Range("L138").Value = 23609.328
$ws.Range("M138").Value = -3961.499899999999
$ws.Range("N138").Value = -33889.328

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18854.75
$ws.Range("I32").Value = 15898.154
$ws.Range("K32").Value = 15898.154
$ws.Range("M32").Value = -15611.154
$ws.Range("H122").Value = 3355.0908
$ws.Range("I122").Value = 2134.6428
$ws.Range("J122").Value = 5490.875
$ws.Range("K122").Value = 6403.928400000001
$ws.Range("L122").Value = 16472.625
$ws.Range("M122").Value = -3953.928400000001
$ws.Range("N122").Value = -21372.625
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents() | Out-Null

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents() | Out-Null
$ws.Range("H99").Value = 5061.9287
$ws.Range("I99").Value = 3785.7
$ws.Range("K99").Value = 3785.7
$ws.Range("M99").Value = -2287.7
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents() | Out-Null
$ws.Range("H134").Value = 2672.0417
$ws.Range("I134").Value = 2525.75
$ws.Range("J134").Value = 4281.25
$ws.Range("K134").Value = 7577.25
$ws.Range("L134").Value = 12843.75
$ws.Range("M134").Value = -5042.25
$ws.Range("N134").Value = -17913.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1133.3334
$ws.Range("I15").Value = 1000
$ws.Range("K15").Value = 1000
$ws.Range("M15").Value = -830
$ws.Range("H22").Value = 1124.4706
$ws.Range("I22").Value = 546.9091
$ws.Range("K22").Value = 546.9091
$ws.Range("M22").Value = -196.9091
$ws.Range("H31").Value = 2504099.5
$ws.Range("I31").Value = 3848489.2
$ws.Range("J31").Value = 7376
$ws.Range("K31").Value = 3848489.2
$ws.Range("L31").Value = 7376
$ws.Range("M31").Value = -3848194.2
$ws.Range("N31").Value = -7966
$ws.Range("H34").Value = 2504099.5
$ws.Range("I34").Value = 3848489.2
$ws.Range("J34").Value = 7376
$ws.Range("K34").Value = 3848489.2
$ws.Range("L34").Value = 7376
$ws.Range("M34").Value = -3848287.2
$ws.Range("N34").Value = -7780
$ws.Range("H93").Value = 18703.5
$ws.Range("I93").Value = 18703.5
$ws.Range("K93").Value = 18703.5
$ws.Range("M93").Value = -16831.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2333.3333
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = -2831
$ws.Range("N25").Value = -9338
$ws.Range("H30").Value = 2333.3333
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 3000
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 9000
$ws.Range("M30").Value = -2898
$ws.Range("N30").Value = -9204
$ws.Range("H34").Value = 6355.9473
$ws.Range("I34").Value = 66.666664
$ws.Range("J34").Value = 9258.691999999999
$ws.Range("K34").Value = 199.999992
$ws.Range("L34").Value = 27776.076
$ws.Range("M34").Value = -115.999992
$ws.Range("N34").Value = -27944.076
$ws.Range("H39").Value = 2796
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2796
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 8388
$ws.Range("M39").ClearContents() | Out-Null
$ws.Range("N39").Value = -8976
$ws.Range("H55").Value = 1222.2222
$ws.Range("I55").Value = 683.3333
$ws.Range("J55").Value = 2300
$ws.Range("K55").Value = 2049.9999
$ws.Range("L55").Value = 6900
$ws.Range("M55").Value = -1872.9999
$ws.Range("N55").Value = -7254
$ws.Range("H75").Value = 1277.9231
$ws.Range("I75").Value = 514
$ws.Range("J75").Value = 2169.1667
$ws.Range("K75").Value = 1542
$ws.Range("L75").Value = 6507.500100000001
$ws.Range("M75").Value = -544
$ws.Range("N75").Value = -8503.500100000001
$ws.Range("H78").Value = 1277.9231
$ws.Range("I78").Value = 514
$ws.Range("J78").Value = 2169.1667
$ws.Range("K78").Value = 4626
$ws.Range("L78").Value = 19522.5003
$ws.Range("M78").Value = 366
$ws.Range("N78").Value = -29506.5003
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents() | Out-Null
$ws.Range("H102").Value = 200
$ws.Range("I102").Value = 200
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 600
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 1834
$ws.Range("N102").ClearContents() | Out-Null
$ws.Range("H104").Value = 2985.3333
$ws.Range("J104").Value = 2985.3333
$ws.Range("L104").Value = 8955.999899999999
$ws.Range("N104").Value = -14197.9999
$ws.Range("H137").Value = 2641.36
$ws.Range("I137").Value = 2473.6316
$ws.Range("J137").Value = 3172.5
$ws.Range("K137").Value = 7420.8948
$ws.Range("L137").Value = 9517.5
$ws.Range("M137").Value = -2320.8948
$ws.Range("N137").Value = -19717.5
$ws.Range("H139").Value = 13892873
$ws.Range("I139").Value = 14709101
$ws.Range("K139").Value = 44127303
$ws.Range("M139").Value = -44122163
$ws.Range("H140").Value = 83334340
$ws.Range("I140").Value = 83334340
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 250003020
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -249997840
$ws.Range("N140").ClearContents() | Out-Null
$ws.Range("H141").Value = 2333.3333
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 31173.75
$ws.Range("I102").Value = 2839.3845
$ws.Range("J102").Value = 104843.1
$ws.Range("K102").Value = 2839.3845
$ws.Range("L102").Value = 104843.1
$ws.Range("M102").Value = -1217.3845
$ws.Range("N102").Value = -108087.1
$ws.Range("H132").Value = 4634.8237
$ws.Range("J132").Value = 4064.8333
$ws.Range("L132").Value = 12194.4999
$ws.Range("N132").Value = -17254.4999
$ws.Range("H133").Value = 29356.666
$ws.Range("J133").Value = 29356.666
$ws.Range("L133").Value = 29356.666
$ws.Range("N133").Value = -39476.666
$ws.Range("H135").Value = 29746.666
$ws.Range("J135").Value = 29746.666
$ws.Range("L135").Value = 29746.666
$ws.Range("N135").Value = -39886.666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4170.125
$ws.Range("I122").Value = 2892.75
$ws.Range("J122").Value = 5447.5
$ws.Range("K122").Value = 8678.25
$ws.Range("L122").Value = 16342.5
$ws.Range("M122").Value = -6228.25
$ws.Range("N122").Value = -21242.5
$ws.Range("H132").Value = 4609.421
$ws.Range("I132").Value = 3676
$ws.Range("J132").Value = 4942.7856
$ws.Range("K132").Value = 11028
$ws.Range("L132").Value = 14828.3568
$ws.Range("M132").Value = -8498
$ws.Range("N132").Value = -19888.3568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 32420.715
$ws.Range("J138").Value = 32420.715
$ws.Range("L138").Value = 32420.715
$ws.Range("N138").Value = -42700.715
